$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.919.59'
$ws.Range("E2").Value = '  -0.55%  '
$ws.Range("D3").Value = '2.040.39'
$ws.Range("E3").Value = '  -0.24%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '247.57'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.70%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.660'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.36%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '55.94'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.74%  '
$ws.Range("E9").Value = '  -0.85%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0776'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +2.18%  '
$ws.Range("E11").Value = '  +1.37%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.62'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +3.07%  '
$ws.Range("D13").Value = '2.339.83'
$ws.Range("E13").Value = '  -0.18%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.58'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +5.55%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.785'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -5.08%  '
$ws.Range("D16").Value = '2.042.61'
$ws.Range("E16").Value = '  -0.18%  '
$ws.Range("D17").Value = '36.906.08'
$ws.Range("E17").Value = '  -0.48%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '16.29'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +12.76%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '73.51'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.31%  '
$ws.Range("D20").Value = '0.0₃0890'
$ws.Range("E20").Value = '  +1.62%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.29'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '235.21'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.35%  '
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.34'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -3.80%  '
$ws.Range("E25").Value = '  +7.58%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '167.31'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.96%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.02'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -1.86%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.70'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -3.31%  '
$ws.Range("E29").Value = '  +0.38%  '
$ws.Range("E30").Value = '  +2.78%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.63'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.98%  '
$ws.Range("E32").Value = '  -3.49%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.38'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.27%  '
$ws.Range("E34").Value = '  -0.25%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0867'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +1.68%  '
$ws.Range("E36").Value = '  -3.67%  '
$ws.Range("E37").Value = '  -1.69%  '
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.33'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.73%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.21'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +15.07%  '
$ws.Range("B40").Value = 'Cronos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.105'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -4.60%  '
$ws.Range("E41").Value = '  +22.14%  '
$ws.Range("E42").Value = '  -3.14%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '17.08'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -6.59%  '
$ws.Range("B44").Value = 'ARBITRUM'
$ws.Range("C44").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.11'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -3.45%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '94.83'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -2.72%  '
$ws.Range("E46").Value = '  +1.07%  '
$ws.Range("D47").Value = '1.274.62'
$ws.Range("E47").Value = '  -2.53%  '
$ws.Range("E48").Value = '  -2.13%  '
$ws.Range("D49").Value = '2.226.19'
$ws.Range("E49").Value = '  -0.52%  '
$ws.Range("E50").Value = '  -3.65%  '
$ws.Range("B51").Value = 'FTXToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.34'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -23.13%  '
